$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet.
#    Seed it from the "2021-Q4" sheet so it inherits the same header style,
#    A-column index style, and per-column formatting (text columns for
#    B-G, numeric column for H), then overwrite the cell values.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totals)
$q1.Name = "2022-Q1"

$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("B1:H1").Copy($q1.Range("B1:H1"))
$template.Range("A2:H5").Copy($q1.Range("A2:H5"))

# Row 2 - 009414 中银大健康股票A
$q1.Cells.Item(2,2).Value = "'009414"
$q1.Cells.Item(2,2).ClearFormats()
$q1.Cells.Item(2,3).Value = "中银大健康股票A"
$q1.Cells.Item(2,4).Value = "'2.96"
$q1.Cells.Item(2,4).ClearFormats()
$q1.Cells.Item(2,5).Value = "'86.25"
$q1.Cells.Item(2,5).ClearFormats()
$q1.Cells.Item(2,6).Value = "'2.43"
$q1.Cells.Item(2,6).ClearFormats()
$q1.Cells.Item(2,7).Value = "'0.0719"
$q1.Cells.Item(2,7).ClearFormats()
$q1.Cells.Item(2,8).Value = 9

# Row 3 - 010663 长江均衡成长混合A
$q1.Cells.Item(3,2).Value = "'010663"
$q1.Cells.Item(3,2).ClearFormats()
$q1.Cells.Item(3,3).Value = "长江均衡成长混合A"
$q1.Cells.Item(3,4).Value = "'0.26"
$q1.Cells.Item(3,4).ClearFormats()
$q1.Cells.Item(3,5).Value = "'85.90"
$q1.Cells.Item(3,5).ClearFormats()
$q1.Cells.Item(3,6).Value = "'4.63"
$q1.Cells.Item(3,6).ClearFormats()
$q1.Cells.Item(3,7).Value = "'0.0120"
$q1.Cells.Item(3,7).ClearFormats()
$q1.Cells.Item(3,8).Value = 3

# Row 4 - 010321 中银大健康股票C
$q1.Cells.Item(4,2).Value = "'010321"
$q1.Cells.Item(4,2).ClearFormats()
$q1.Cells.Item(4,3).Value = "中银大健康股票C"
$q1.Cells.Item(4,4).Value = "'0.10"
$q1.Cells.Item(4,4).ClearFormats()
$q1.Cells.Item(4,5).Value = "'86.25"
$q1.Cells.Item(4,5).ClearFormats()
$q1.Cells.Item(4,6).Value = "'2.43"
$q1.Cells.Item(4,6).ClearFormats()
$q1.Cells.Item(4,7).Value = "'0.0024"
$q1.Cells.Item(4,7).ClearFormats()
$q1.Cells.Item(4,8).Value = 9

# Row 5 - 010664 长江均衡成长混合C
$q1.Cells.Item(5,2).Value = "'010664"
$q1.Cells.Item(5,2).ClearFormats()
$q1.Cells.Item(5,3).Value = "长江均衡成长混合C"
$q1.Cells.Item(5,4).Value = "'0.05"
$q1.Cells.Item(5,4).ClearFormats()
$q1.Cells.Item(5,5).Value = "'85.90"
$q1.Cells.Item(5,5).ClearFormats()
$q1.Cells.Item(5,6).Value = "'4.63"
$q1.Cells.Item(5,6).ClearFormats()
$q1.Cells.Item(5,7).Value = "'0.0023"
$q1.Cells.Item(5,7).ClearFormats()
$q1.Cells.Item(5,8).Value = 3

# ---------------------------------------------------------------------------
# 2) Add the corresponding 2022-Q1 summary row to the "总计" sheet: insert
#    a new row above the existing 2021-Q4 row, copy the index-cell style
#    from the row below it, fill in the new data, then bump the existing
#    running index numbers in column A by one.
#    NOTE: re-fetch the "总计" worksheet reference here - the handle
#    obtained before Worksheets.Add()/rename can end up bound to whichever
#    sheet is now active instead of staying pinned to the original sheet.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()
$totals.Cells.Item(3,1).Copy($totals.Cells.Item(2,1))
$totals.Range("B2:D2").ClearFormats()

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q1"
$totals.Cells.Item(2,3).Value = 4
$totals.Cells.Item(2,4).Value = 0.09

$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(6,1).Value = 4
